# Weekly driver report update for 2025-04-28
#
# Refreshes the "Bad Drivers" and "Good Drivers" tables on the
# "Driver Summary" sheet with this week's roaming-impact numbers, and
# drops the drivers that aged out of the "Good Drivers" list (the sheet
# shrinks from 26 used rows down to 20).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Driver Summary")

# ---- Bad Drivers table (rows 3-5) ----
# Row 3 adapter name is unchanged, only its stats move.
$ws.Range("C3").Value = 383
$ws.Range("D3").Value = 93.3

# Row 4 adapter was re-identified from AX201 to AX211 this week.
$ws.Range("A4").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 23.30.0.6"
$ws.Range("B4").Value = 50
$ws.Range("C4").Value = 1074
$ws.Range("D4").Value = 98.7

# Row 5 totals
$ws.Range("B5").Value = 53
$ws.Range("C5").Value = 1457

# ---- Good Drivers table (rows 13-20) ----
# Row 13: previous #2 driver moves up to #1, with refreshed counts.
$ws.Range("A13").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.3.1"
$ws.Range("B13").Value = 11140
$ws.Range("D13").Value = 100
$ws.Range("E13").Value = "'2022-08-29"

# Row 14: previous #3 driver moves up to #2.
$ws.Range("A14").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.0.3"
$ws.Range("B14").Value = 14487
$ws.Range("D14").Value = 100
$ws.Range("E14").Value = "'2022-05-23"

# Row 15: previous #4 driver moves up to #3.
$ws.Range("A15").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.100.1.1"
$ws.Range("B15").Value = 265400
$ws.Range("D15").Value = 99.90000000000001
$ws.Range("E15").Value = "'2022-05-01"

# The remaining six "Good Drivers" entries (previously rows 16-21) aged
# off the report this week, and the trailing blank rows 22-26 go with
# them, so the sheet's used range shrinks to row 20.
$ws.Range("A16:J26").Clear()
